# Adds a "fail"/"pass" verdict column (R) to the "Input Data" sheet and
# updates three e-mail addresses (and their hyperlinks) in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input Data")

# --- New column R: give it the same "Text" number format used elsewhere
# in the sheet (column D etc. use style index 1 / numFmtId 49) ---
$ws.Range("R1:R7").NumberFormat = "@"

# Fill column R with pass/fail verdicts, one per data row.
# "R3" is written first so that "fail" lands in the shared-string table
# ahead of "pass" (matches original authoring order).
$ws.Range("R3").Value = "fail"
$ws.Range("R1").Value = "pass"
$ws.Range("R2").Value = "pass"
$ws.Range("R4").Value = "fail"
$ws.Range("R5").Value = "fail"
$ws.Range("R6").Value = "pass"
$ws.Range("R7").Value = "fail"

# --- Update the three e-mail addresses in column D that changed ---
$ws.Range("D1").Value = "email@ss.com52"
$ws.Range("D2").Value = "email@ss.com6455"
$ws.Range("D6").Value = "email@ss.com101111"

# New explicit hyperlinks were added on top of the existing D2:D3 / D5:D6
# hyperlink ranges for the two cells whose address text changed.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:email@ss.com6455")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:email@ss.com101111")

# Adding a hyperlink re-styles the cell; put the Hyperlink cell style back
# so D2/D6 keep looking like the rest of the (already-hyperlinked) column.
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"

# Match the saved cursor position recorded in the workbook.
$ws.Range("E13").Select()
